$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("en")

# Append the new key/value pairs to the language table, following the
# existing pattern (column A = key, column B = value with wrap-text style).
$ws.Range("B21").Value = "Multiply Error:"
$ws.Range("B21").WrapText = $true
$ws.Range("A21").Value = "victory_errorMult"

$ws.Range("A22").Value = "victory_errorSums"
$ws.Range("B22").Value = "Sums Error:"
$ws.Range("B22").WrapText = $true

$ws.Range("A23").Value = "victory_score"
$ws.Range("B23").Value = "Score:"
$ws.Range("B23").WrapText = $true

# Move the active selection to match the new last row, mirroring the
# workbook's recorded UI state after the edit.
$ws.Range("A23").Select()
